$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Populate new row 43 with the "Ditto / EFT Generation" entry ---
$dt = Get-Date -Year 2015 -Month 4 -Day 26 -Hour 0 -Minute 0 -Second 0
$ws.Range("A43").Value = $dt
$ws.Range("A43").NumberFormat = "m/d/yy"

$ws.Range("C43").Value = "New Files`nEJBServer\components\MOLSA\message\MOLSABpoNumberToArabic.xml`nEJBServer\components\MOLSA\source\curam\molsa\util\impl\MOLSANumberToArabic.java`nChanged Files`nEJBServer\components\MOLSA\message\MOLSABpoGenerateEFT.xml`nEJBServer\components\MOLSA\source\curam\molsa\util\impl\MOLSAGenerateEFTHelper.java`nEJBServer\components\MOLSA\source\curam\molsa\eft\batch\impl\MOLSAGenerateEFTBatchStream.java`nEJBServer\components\MOLSA\properties\Application.prx`nmodel/Packages/EFT/Batch.efx"

$ws.Range("D43").Value = "EFT Generation"
$ws.Range("E43").Value = "Ditto"
$ws.Range("F43").Value = "Yes"
$ws.Range("G43").Value = "Only the Application.prx change (use insertproperties and change back the production details)"
$ws.Range("I43").Value = "Use build insertProperties and change back the production varaiables."

# Row 43 needs extra height to show the wrapped multi-line comment
$ws.Rows("43").RowHeight = 195

# --- Update the view so the freeze pane / selection matches the new layout ---
$ws.Range("C43").Select()
